# TMTT0011424 test-data refresh (Hierarchy Viewer + Time recordManager):
#  - Users: rename the sample associate used by the TAS supervisor test
#  - Project_Title: preserve the previous dropdown value in a new "Old
#    Project Dropdown Selection" column and roll the sample project forward
#  - RateSheetManagement / WeeklyEntryMatrix: rename the sample engagement
#    to the new LucidHealth project code

$wb = $excel.ActiveWorkbook

# --- Users ---------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Activate()
$wsUsers.Range("A2").Value = "Christy Skaar"
$wsUsers.Columns.Item(1).AutoFit()
[void]$wsUsers.Range("A2").Select()

# --- Project_Title --------------------------------------------------------
$wsProject = $wb.Worksheets.Item("Project_Title")
$wsProject.Activate()
$wsProject.Range("C1").Value = "Old Project Dropdown Selection"
$wsProject.Range("C1").Font.Bold = $true
$wsProject.Range("C2").Value = "Engagement Project Bend-FVA-109081"
$wsProject.Range("A2").Value = "Project Newman-United Flow Technologies, LLC-FVA-115826"
$wsProject.Columns.Item(1).AutoFit()
$wsProject.Columns.Item(2).AutoFit()
$wsProject.Columns.Item(3).AutoFit()
[void]$wsProject.Range("C2").Select()

# --- StaffMember (visited, no data change) --------------------------------
$wsStaff = $wb.Worksheets.Item("StaffMember")
$wsStaff.Activate()

# --- Update_Hours (visited, no data change) -------------------------------
$wsHours = $wb.Worksheets.Item("Update_Hours")
$wsHours.Activate()

# --- Update_Timer (visited, no data change; loses the old active tab) ----
$wsTimer = $wb.Worksheets.Item("Update_Timer")
$wsTimer.Activate()
[void]$wsTimer.Range("A2").Select()

# --- RateSheetManagement ---------------------------------------------------
$wsRateSheet = $wb.Worksheets.Item("RateSheetManagement")
$wsRateSheet.Activate()
$wsRateSheet.Range("A2").Value = "Project Clear-LucidHealth-FVA-105379"
$wsRateSheet.Columns.Item(1).AutoFit()
[void]$wsRateSheet.Range("A2").Select()

# --- WeeklyEntryMatrix (final active sheet) --------------------------------
$wsWeekly = $wb.Worksheets.Item("WeeklyEntryMatrix")
$wsWeekly.Activate()
$wsWeekly.Range("A2").Value = "Project Clear-LucidHealth-FVA-105379"
$wsWeekly.Columns.Item(1).AutoFit()
[void]$wsWeekly.Range("A2").Select()
